# "This is modified my excel"
# Populate Sheet1!C5:E5 with the name fields and leave E5 selected,
# matching the author's edit (Name: / Shweta  / Pitambare).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = "Name:"
$ws.Range("D5").Value = "Shweta "
$ws.Range("E5").Value = "Pitambare"

$ws.Range("E5").Select()
